$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packs")
$ws.Activate()

# New row 34 - "promoted1" offer pack definition (mirrors the "promoted3" row but
# as a featured/enabled-by-default promo, with a new min app version of 1.12).
$ws.Range("B34").Value = "<Definition>"
$ws.Range("C34").Value = "promoted1"
$ws.Range("D34").Value = "promoted1"
$ws.Range("E34").Value = $true
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = $false
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = 1.99
$ws.Range("S34").Value = 0.6
$ws.Range("T34").Value = "com.ubisoft.hungrydragon.promoted1"
$ws.Range("U34").Value = "TID_OFFER_PACK_SUPER"
$ws.Range("V34").Value = $true
$ws.Range("W34").Value = 0
$ws.Range("X34").Value = "dragonSelection"
$ws.Range("Y34").Value = 10
$ws.Range("Z34").Value = "-"
$ws.Range("AA34").Value = "-"
$ws.Range("AB34").Value = "-"
$ws.Range("AC34").Value = "1.12"
$ws.Range("AD34").Value = "-"
$ws.Range("AE34").Value = "-"
$ws.Range("AF34").Value = 4
$ws.Range("AG34").Value = "-"
$ws.Range("AH34").Value = "-"
$ws.Range("AI34").Value = "-"
$ws.Range("AJ34").Value = "-"
$ws.Range("AK34").Value = "-"
$ws.Range("AL34").Value = "-"
$ws.Range("AM34").Value = "-"
$ws.Range("AN34").Value = "-"
$ws.Range("AO34").Value = "-"
$ws.Range("AP34").Value = "-"
$ws.Range("AQ34").Value = "-"
$ws.Range("AR34").Value = "-"
$ws.Range("AS34").Value = "-"
$ws.Range("AT34").Value = "-"
$ws.Range("AU34").Value = "-"
$ws.Range("AV34").Value = "-"
$ws.Range("AW34").Value = "-"

# Grow the table to include the new row.
$lo = $ws.ListObjects.Item("offerPacksDefinitions")
$lo.Resize($ws.Range("B2:AW34"))

# Restore the view position/selection recorded for the sheet.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("H35").Select()
